$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the nurse's record (row 2): name, phone, password, and status.
# Leading apostrophes force Excel to store the numeric-looking phone
# number and password as text (matching the original shared-string type)
# instead of silently converting them to numbers.
$ws.Range("B2").Value = "Abdullayeva Maftuna"
$ws.Range("C2").Value = "'+998971234566"
$ws.Range("D2").Value = "'1234"
$ws.Range("F2").Value = "ACTIVE"
